# Updated cryptos list on Mon Dec  4 03:25:41 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.794.84"
$ws.Range("E2").Value = "  +3.67%  "
$ws.Range("D3").Value = "2.216.06"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'229.20"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("D7").Value = "'64.54"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.406"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").Value = "'0.0870"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "2.545.19"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "'15.92"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "'22.26"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "'0.821"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "2.206.97"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "40.686.39"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "'74.37"
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  +6.01%  "
$ws.Range("D21").Value = "'6.15"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'254.47"
$ws.Range("E22").Value = "  +10.04%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'2.37"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  -8.22%  "
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("D27").Value = "'173.31"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("D29").Value = "'20.40"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("D31").Value = "'2.81"
$ws.Range("E31").Value = "  +3.36%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").Value = "'7.17"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").Value = "'4.79"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  +6.83%  "
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'4.91"
$ws.Range("E40").Value = "  +13.82%  "
$ws.Range("D41").Value = "'8.67"
$ws.Range("E41").Value = "  +11.16%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "'101.23"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "'1.24"
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("D45").Value = "1.521.11"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "'17.34"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").Value = "'0.0938"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.82"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").Value = "'0.000206"
$ws.Range("E50").Value = "  +38.75%  "
$ws.Range("B51").Value = "Celestia"
$ws.Range("C51").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D51").Value = "'9.54"
$ws.Range("E51").Value = "  +10.37%  "